$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cellRef, $value) {
    $rng = $ws.Range($cellRef)
    $rng.NumberFormat = "@"
    $rng.Value = $value
    $rng.ClearFormats()
}

Set-TextValue "D2" "274.69"
Set-TextValue "E2" "2.31%"
Set-TextValue "D3" "26.81"
Set-TextValue "E3" "0.55%"
Set-TextValue "D4" "4.916"
Set-TextValue "E4" "4.50%"
Set-TextValue "D5" "0.06358"
Set-TextValue "E5" "4.28%"
Set-TextValue "D6" "6.961"
Set-TextValue "E6" "3.44%"
$ws.Range("B7").Value = "FTXToken"
$ws.Range("C7").Value = "https://coinranking.com/coin/NfeOYfNcl+ftxtoken-ftt"
Set-TextValue "D7" "1.502"
Set-TextValue "E7" "67.74%"
$ws.Range("B8").Value = "GateToken"
$ws.Range("C8").Value = "https://coinranking.com/coin/t7m8DZVyMsAu+gatetoken-gt"
Set-TextValue "D8" "3.357"
Set-TextValue "E8" "5.99%"
Set-TextValue "D9" "0.8877"
Set-TextValue "E9" "3.58%"
Set-TextValue "D10" "0.1473"
Set-TextValue "E10" "4.09%"
Set-TextValue "D11" "0.05209"
Set-TextValue "E11" "6.22%"
Set-TextValue "D12" "0.07415"
Set-TextValue "E12" "4.45%"
Set-TextValue "D13" "0.03163"
Set-TextValue "E13" "-1.71%"
Set-TextValue "D14" "0.09056"
Set-TextValue "E14" "0.43%"
Set-TextValue "D15" "0.001564"
Set-TextValue "E15" "2.27%"
Set-TextValue "D16" "0.0006326"
Set-TextValue "E16" "4.21%"
Set-TextValue "D17" "0.006027"
Set-TextValue "E17" "0.02%"
Set-TextValue "D18" "3.485"
Set-TextValue "E18" "0.69%"
Set-TextValue "E19" "1.72%"
Set-TextValue "E20" "2.23%"
Set-TextValue "D21" "0.1334"
Set-TextValue "E21" "2.74%"
Set-TextValue "D22" "3.953"
Set-TextValue "E22" "2.79%"
Set-TextValue "D23" "0.04340"
Set-TextValue "E23" "2.70%"
Set-TextValue "D24" "0.001185"
Set-TextValue "E24" "0.36%"
Set-TextValue "D25" "0.003663"
Set-TextValue "E25" "-11.63%"
Set-TextValue "D26" "0.0001204"
Set-TextValue "E26" "0.37%"
Set-TextValue "D27" "0.0001942"
Set-TextValue "E27" "15.66%"
Set-TextValue "D40" "0.04038"
Set-TextValue "E40" "2.41%"
Set-TextValue "D41" "0.006640"
Set-TextValue "E41" "58.71%"
Set-TextValue "D42" "0.1169"
Set-TextValue "E42" "4.70%"
Set-TextValue "D43" "0.002369"
Set-TextValue "E43" "17.90%"
Set-TextValue "D44" "0.01226"
Set-TextValue "E44" "-2.80%"
Set-TextValue "D45" "0.00005252"
Set-TextValue "E46" "442.58%"
Set-TextValue "D47" "0.02126"
Set-TextValue "E47" "-13.13%"
Set-TextValue "E48" "0.02%"
